$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# ---------------------------------------------------------------------------
# Build the two new border styles once on sheet1 (quality_comparison):
#   C1 -> thin top + thin bottom               (matches new borderId 4)
#   D1 -> thin top + thin bottom + thin right   (matches new borderId 5)
# ---------------------------------------------------------------------------
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$d1.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$d1.Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# Re-use those exact styles (instead of re-deriving the borders edge-by-edge,
# which the engine would otherwise record as brand-new, if equivalent, style
# records) by copying the formats onto the remaining header cells.
$ws1.Range("C1").Copy() | Out-Null
$ws2.Range("C1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122) | Out-Null

$ws1.Range("D1").Copy() | Out-Null
$ws2.Range("D1").PasteSpecial(-4122) | Out-Null
$ws2.Range("G1").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Anonymize "fedcore" -> "approach"
# ---------------------------------------------------------------------------
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# ---------------------------------------------------------------------------
# Drop the stray empty inline-string cell at G5 on computational_comparison
# ---------------------------------------------------------------------------
$ws2.Range("G5").ClearContents()

Write-Host "edits applied"
